$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) & 2) UVOD section: remove the "Tema ovog projekta..." and
#    "U nasem slucaju..." paragraphs, and replace the following
#    (previously empty, page-break-only) paragraph's content with
#    "(...)" split around a relocated "_GoBack" bookmark. Adding a
#    bookmark with the same name elsewhere automatically drops the
#    old "_GoBack" bookmark located a few paragraphs above
#    ("24. 0[_GoBack]3. 2019."), taking care of that removal too.
# -----------------------------------------------------------------

$headingRange = $d.Content
$headingRange.Find.Execute("UVOD", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$searchRange = $d.Range($headingRange.End, $d.Content.End)
$searchRange.Find.Execute("Tema ovog projekta je Association Rule Mining", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

# Expand to the full first paragraph ("Tema ovog projekta...").
$searchRange.Expand(4) | Out-Null
$firstParaEnd = $searchRange.End

# Expand to the full second paragraph ("U nasem slucaju...").
$secondPara = $d.Range($firstParaEnd, $firstParaEnd)
$secondPara.Expand(4) | Out-Null

# Delete both paragraphs in one shot.
$bothParas = $d.Range($searchRange.Start, $secondPara.End)
$targetStart = $bothParas.Start
$bothParas.Delete()

# $targetStart now sits at the start of the paragraph that used to
# hold only a page break. Insert ")" then "(..." before it so the
# final text reads "(...)" immediately before the page break.
$closeParen = $d.Range($targetStart, $targetStart)
$closeParen.InsertBefore(")")

$openParen = $d.Range($targetStart, $targetStart)
$openParen.InsertBefore("(...")

# Re-anchor the (hidden) "_GoBack" bookmark between "(..." and ")".
$bmPos = $targetStart + 4
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# -----------------------------------------------------------------
# 3) "Pouzdanost" section: collapse the many small runs making up
#    the "Softver smatramo pouzdanim..." paragraph into one run
#    with the full text (unchanged wording, just de-fragmented).
# -----------------------------------------------------------------

$relRange = $d.Content
$relRange.Find.Execute("oftver smatramo pouzdanim ukoliko", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$relRange.MoveStart(1, -1) | Out-Null   # back up over the leading "S"
$relRange.Expand(4) | Out-Null          # whole paragraph, incl. mark

$paraStart = $relRange.Start
$paraEndNoMark = $relRange.End - 1

$fullText = "Softver smatramo pouzdanim ukoliko obavlja svoju zada" + [char]0x0107 + "u u skladu sa specifikacijom bez gre" + [char]0x0161 + "aka koje mogu dovesti do ne" + [char]0x017E + "eljenih rezultata ili potpunog prestanka rada sustava. U praksi softversku pouzdanost " + [char]0x010D + "e" + [char]0x0161 + [char]0x0107 + "e izra" + [char]0x017E + "avamo kao postotak ili o" + [char]0x010D + "ekivano vrijeme rada bez gre" + [char]0x0161 + "aka."

$paraRange = $d.Range($paraStart, $paraEndNoMark)
$paraRange.Text = $fullText
